# Insert a new weekly price record for "Plátano" at row 366, shifting the
# existing rows 366-445 down to 367-446 (dimension grows from A1:T445 to
# A1:T446). The new row carries the same market/product metadata as the
# surrounding rows but a new date and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 366..445 down to 367..446, opening up a blank row 366.
$ws.Rows("366:366").Insert()

# Populate the newly-opened row 366 with the new weekly record.
$ws.Cells.Item(366, 1).Value2  = 1
$ws.Cells.Item(366, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(366, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(366, 4).Value2  = 45244
$ws.Cells.Item(366, 5).Value2  = 15
$ws.Cells.Item(366, 6).Value2  = "Fruta"
$ws.Cells.Item(366, 7).Value2  = 100108
$ws.Cells.Item(366, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(366, 9).Value2  = 100108006
$ws.Cells.Item(366, 10).Value2 = "Plátano"
$ws.Cells.Item(366, 11).Value2 = "Sin especificar"
$ws.Cells.Item(366, 12).Value2 = "Pintón"
$ws.Cells.Item(366, 13).Value2 = 120
$ws.Cells.Item(366, 14).Value2 = 30000
$ws.Cells.Item(366, 15).Value2 = 32000
$ws.Cells.Item(366, 16).Value2 = 31000
$ws.Cells.Item(366, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(366, 18).Value2 = "Ecuador"
$ws.Cells.Item(366, 19).Value2 = 1550
$ws.Cells.Item(366, 20).Value2 = 20
